# Update the "Test" table from a 3-column x 2-row stub into the 3x9
# order table: rename the original columns, add 6 more generic
# "Column" columns, and append 6 more order rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Give the first three columns their real headers.
$ws.Range("A1").Value = "First Name"
$ws.Range("B1").Value = "Last Name"
$ws.Range("C1").Value = "Special #"

# Grow the table by 6 columns (D..I) -> 9 columns total.
for ($i = 0; $i -lt 6; $i++) {
  $tbl.ListColumns.Add() | Out-Null
}

# The newly added columns keep the tool's generic "Column" name.
$ws.Range("D1:I1").Value = "Column"

# Append rows 4..9 of order data. Copy the existing data row down so the
# "0" placeholder stays text (shared string) like the rest of the sheet,
# instead of being re-interpreted as a number.
for ($r = 4; $r -le 9; $r++) {
  $ws.Range("A2:C2").Copy($ws.Range("A" + $r + ":C" + $r))
}
